$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '70.121.31'
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").Value = '2.473.59'
$ws.Range("E3").Value = '  +0.86%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '''569.27'
$ws.Range("E5").Value = '  +2.08%  '
$ws.Range("D6").Value = '''167.53'
$ws.Range("E6").Value = '  +3.05%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '''0.514'
$ws.Range("E8").Value = '  +0.83%  '
$ws.Range("D9").Value = '''0.178'
$ws.Range("E9").Value = '  +13.90%  '
$ws.Range("E10").Value = '  -1.27%  '
$ws.Range("D11").Value = '''0.336'
$ws.Range("E11").Value = '  +2.45%  '
$ws.Range("D12").Value = '''4.69'
$ws.Range("E12").Value = '  -2.65%  '
$ws.Range("D13").Value = '''0.0000185'
$ws.Range("E13").Value = '  +9.20%  '
$ws.Range("D14").Value = '70.004.96'
$ws.Range("E14").Value = '  +2.55%  '
$ws.Range("D15").Value = '2.921.00'
$ws.Range("E15").Value = '  +0.52%  '
$ws.Range("D16").Value = '''24.12'
$ws.Range("E16").Value = '  +3.35%  '
$ws.Range("D17").Value = '2.466.95'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = '''10.87'
$ws.Range("E18").Value = '  +3.72%  '
$ws.Range("D19").Value = '''344.21'
$ws.Range("E19").Value = '  +2.24%  '
$ws.Range("D20").Value = '''7.20'
$ws.Range("E20").Value = '  +4.66%  '
$ws.Range("B21").Value = 'SuiNetwork'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D21").Value = '''2.05'
$ws.Range("E21").Value = '  +9.91%  '
$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").Value = '''3.91'
$ws.Range("E22").Value = '  +3.37%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '''66.71'
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").Value = '''3.92'
$ws.Range("E25").Value = '  +6.86%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").Value = '''8.58'
$ws.Range("E26").Value = '  +5.68%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '2.588.81'
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("D28").Value = '''0.986'
$ws.Range("E28").Value = '  -1.18%  '
$ws.Range("D29").Value = '0.0₃0860'
$ws.Range("E29").Value = '  +5.59%  '
$ws.Range("D30").Value = '''7.38'
$ws.Range("E30").Value = '  +2.65%  '
$ws.Range("D31").Value = '''1.27'
$ws.Range("E31").Value = '  +11.05%  '
$ws.Range("D32").Value = '''452.32'
$ws.Range("E32").Value = '  +6.28%  '
$ws.Range("D33").Value = '''0.998'
$ws.Range("E33").Value = '  -0.22%  '
$ws.Range("D34").Value = '''1.64'
$ws.Range("E34").Value = '  +1.72%  '
$ws.Range("D35").Value = '''160.82'
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("E37").Value = '  +4.38%  '
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").Value = '''18.24'
$ws.Range("E39").Value = '  +2.59%  '
$ws.Range("D40").Value = '''0.308'
$ws.Range("E40").Value = '  +4.39%  '
$ws.Range("D41").Value = '''4.53'
$ws.Range("E41").Value = '  +3.51%  '
$ws.Range("D42").Value = '''1.56'
$ws.Range("E42").Value = '  +6.40%  '
$ws.Range("E43").Value = '  +4.06%  '
$ws.Range("D44").Value = '''2.17'
$ws.Range("E44").Value = '  +7.95%  '
$ws.Range("D45").Value = '''3.42'
$ws.Range("E45").Value = '  +1.82%  '
$ws.Range("D46").Value = '''133.17'
$ws.Range("E46").Value = '  +2.71%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '''0.494'
$ws.Range("E47").Value = '  +2.80%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.0725'
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("D49").Value = '''0.566'
$ws.Range("E49").Value = '  +0.86%  '
$ws.Range("D50").Value = '''0.0932'
$ws.Range("E50").Value = '  +1.59%  '
$ws.Range("E51").Value = '  +2.79%  '
